# Update "想去人数" (F column) values on the "展览" sheet and the
# "全部类型" sheet, which both contain the same event listing.

$wb = $excel.ActiveWorkbook

# Values that are identical across both updated sheets.
$common = @{
    2  = 136
    3  = 1658
    4  = 642
    5  = 1111
    6  = 22
    7  = 11670
    8  = 34
    9  = 95
    10 = 464
    11 = 382
    15 = 13229
    22 = 46
    23 = 138
}

# "展览" sheet (sheet1) - F14 goes to 13133
$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $common.Keys) {
    $ws1.Cells.Item($row, 6).Value = $common[$row]
}
$ws1.Cells.Item(14, 6).Value = 13133

# "全部类型" sheet (sheet4) - F14 goes to 13139
$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $common.Keys) {
    $ws4.Cells.Item($row, 6).Value = $common[$row]
}
$ws4.Cells.Item(14, 6).Value = 13139
